$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 30,9
$data[0,0] = '2025-07-10'
$data[0,1] = 98
$data[0,2] = 'RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA'
$data[0,3] = 54141318
$data[0,4] = 11939681
$data[0,5] = 'VASSOURA VARRE CANTO COM CABO PLASTIFICADO'
$data[0,6] = 72
$data[0,7] = 5.86
$data[0,8] = 11.18
$data[1,0] = '2025-07-10'
$data[1,1] = 35
$data[1,2] = 'SAT BRAS INDUSTRIA ELETRONICA DA AMAZONIA LTDA.'
$data[1,3] = 54142176
$data[1,4] = 11938367
$data[1,5] = 'PAPEL HIGIENICO 8X300 NEWPAPER 100% Celulose'
$data[1,6] = 66
$data[1,7] = 4.76
$data[1,8] = 5.94
$data[2,0] = '2025-07-10'
$data[2,1] = 95
$data[2,2] = 'SAT BRAS INDUSTRIA ELETRONICA DA AMAZONIA LTDA.'
$data[2,3] = 54142176
$data[2,4] = 19264853
$data[2,5] = 'PAPEL TOALHA INTERFOLHADO 1250 FLS NEWPAPER 100% CELULOSE'
$data[2,6] = 530
$data[2,7] = 8.279999999999999
$data[2,8] = 13.51
$data[3,0] = '2025-07-10'
$data[3,1] = 270
$data[3,2] = 'RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA'
$data[3,3] = 54141318
$data[3,4] = 35118277
$data[3,5] = 'SABAO EM PO ABSOLUTO 400G'
$data[3,6] = 711
$data[3,7] = 18.27
$data[3,8] = 33.82
$data[4,0] = '2025-07-10'
$data[4,1] = 89
$data[4,2] = 'RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA'
$data[4,3] = 54141318
$data[4,4] = 13995639
$data[4,5] = 'ESCOVA OVAL PLASTICA'
$data[4,6] = 58
$data[4,7] = 4.34
$data[4,8] = 10.67
$data[5,0] = '2025-07-10'
$data[5,1] = 30
$data[5,2] = 'CR OBRAS DA CONSTRUCAO LTDA'
$data[5,3] = 54123729
$data[5,4] = 33278408
$data[5,5] = 'AROMATIZANTE LIMPADOR PERF CONC COALA ALGODAO 120ML'
$data[5,6] = 35
$data[5,7] = 5.83
$data[5,8] = 5.79
$data[6,0] = '2025-07-10'
$data[6,1] = 270
$data[6,2] = 'RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA'
$data[6,3] = 54141318
$data[6,4] = 14589837
$data[6,5] = 'FLANELA BRANCA TAM P 28X38CM'
$data[6,6] = 277
$data[6,7] = 15.44
$data[6,8] = 28.27
$data[7,0] = '2025-07-10'
$data[7,1] = 113
$data[7,2] = 'RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA'
$data[7,3] = 54141318
$data[7,4] = 16871438
$data[7,5] = 'DESODORISADOR LADY AEROSSOL 360ML TALCO SUAVE CARINHO'
$data[7,6] = 1452
$data[7,7] = 9.210000000000001
$data[7,8] = 11.36
$data[8,0] = '2025-07-10'
$data[8,1] = 141
$data[8,2] = 'RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA'
$data[8,3] = 54141318
$data[8,4] = 12285275
$data[8,5] = 'LUSTRA MOVEIS BUTTERFLY 200ML AUDAX LAVANDA'
$data[8,6] = 68
$data[8,7] = 9.15
$data[8,8] = 16.36
$data[9,0] = '2025-07-10'
$data[9,1] = 95
$data[9,2] = 'RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA'
$data[9,3] = 54141318
$data[9,4] = 11939551
$data[9,5] = 'RODO COM CABO P 30CM'
$data[9,6] = 50
$data[9,7] = 6.4
$data[9,8] = 15.76
$data[10,0] = '2025-07-10'
$data[10,1] = 112
$data[10,2] = 'RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA'
$data[10,3] = 54141318
$data[10,4] = 16537374
$data[10,5] = 'DESODORISADOR LADY AEROSSOL 360 ML LAVANDA'
$data[10,6] = 2269
$data[10,7] = 8.24
$data[10,8] = 10.09
$data[11,0] = '2025-07-10'
$data[11,1] = 240
$data[11,2] = 'AMMAC INDUSTRIA E COMERCIO DE ALIMENTOS LTDA'
$data[11,3] = 54153624
$data[11,4] = 32130390
$data[11,5] = 'ESPONJA MULTIUSO JEITOSA'
$data[11,6] = 11885
$data[11,7] = 21.36
$data[11,8] = 51.13
$data[12,0] = '2025-07-11'
$data[12,1] = 310
$data[12,2] = 'JURUA ESTALEIROS E NAVEGACAO LTDA'
$data[12,3] = 54124880
$data[12,4] = 12054191
$data[12,5] = 'PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM'
$data[12,6] = 5506
$data[12,7] = 20.1
$data[12,8] = 39.29
$data[13,0] = '2025-07-11'
$data[13,1] = 154
$data[13,2] = 'JURUA ESTALEIROS E NAVEGACAO LTDA'
$data[13,3] = 54124880
$data[13,4] = 11939672
$data[13,5] = 'VASSOURA PIACAVA 20 FUROS'
$data[13,6] = 465
$data[13,7] = 9.99
$data[13,8] = 28.79
$data[14,0] = '2025-07-14'
$data[14,1] = 250
$data[14,2] = 'ERAM ESTALEIRO RIO AMAZONAS LTDA'
$data[14,3] = 54263882
$data[14,4] = 15011531
$data[14,5] = 'DETERGENTE LIMPOL COCO 500ML'
$data[14,6] = 301
$data[14,7] = 20.68
$data[14,8] = 30.82
$data[15,0] = '2025-07-14'
$data[15,1] = 60
$data[15,2] = 'SAWEM DA AMAZONIA LTDA'
$data[15,3] = 54285324
$data[15,4] = 27262762
$data[15,5] = 'MARCA TEXTO AMARELO UND JOCAR OFFICE'
$data[15,6] = 111
$data[15,7] = 7.31
$data[15,8] = 8.16
$data[16,0] = '2025-07-15'
$data[16,1] = 130
$data[16,2] = 'MUSASHI DA AMAZONIA LTDA'
$data[16,3] = 54346779
$data[16,4] = 17171383
$data[16,5] = 'DETERGENTE DESENGRAX MAX PINE AUDAX 5L'
$data[16,6] = -8
$data[16,7] = 14.02
$data[16,8] = 29.5
$data[17,0] = '2025-07-15'
$data[17,1] = 300
$data[17,2] = 'MUSASHI DA AMAZONIA LTDA'
$data[17,3] = 54346779
$data[17,4] = 14795919
$data[17,5] = 'SACO DE LIXO 200L COMUM PACOTINHO C/5 UND SACOLMAX'
$data[17,6] = 244
$data[17,7] = 43.11
$data[17,8] = 61.16
$data[18,0] = '2025-07-16'
$data[18,1] = 70
$data[18,2] = 'TECHLOG - SERVICOS DE GESTAO E SISTEMAS INFORMATIZ'
$data[18,3] = 54316914
$data[18,4] = 19264853
$data[18,5] = 'PAPEL TOALHA INTERFOLHADO 1250 FLS NEWPAPER 100% CELULOSE'
$data[18,6] = 530
$data[18,7] = 8.279999999999999
$data[18,8] = 13.51
$data[19,0] = '2025-07-16'
$data[19,1] = 60
$data[19,2] = 'V V REFEICOES LTDA'
$data[19,3] = 54396269
$data[19,4] = 42173656
$data[19,5] = 'COADOR DE CAFÉ G'
$data[19,6] = 23
$data[19,7] = 6.07
$data[19,8] = 15.64
$data[20,0] = '2025-07-16'
$data[20,1] = 10
$data[20,2] = 'V V REFEICOES LTDA'
$data[20,3] = 54396269
$data[20,4] = 17541022
$data[20,5] = 'COADOR DE CAFE INDUSTRIAL (MAIOR)'
$data[20,6] = 26
$data[20,7] = 2.09
$data[20,8] = 2.05
$data[21,0] = '2025-07-16'
$data[21,1] = 10
$data[21,2] = 'CONDOMINIO DO TVLANDIA MALL'
$data[21,3] = 54119372
$data[21,4] = 28133466
$data[21,5] = 'DESINFETANTE CONCENTRADO 5L AUDAX MAX 1:200 - LAVANDA'
$data[21,6] = 5
$data[21,7] = 2.34
$data[21,8] = 1.91
$data[22,0] = '2025-07-16'
$data[22,1] = 10
$data[22,2] = 'V V REFEICOES LTDA'
$data[22,3] = 54396269
$data[22,4] = 17125814
$data[22,5] = 'COADOR DE CAFE P'
$data[22,6] = 21
$data[22,7] = 2.19
$data[22,8] = 1.98
$data[23,0] = '2025-07-16'
$data[23,1] = 40
$data[23,2] = 'CONDOMINIO RESIDENCIAL EPHYGENIO SALLES'
$data[23,3] = 54365832
$data[23,4] = 20619556
$data[23,5] = 'PEDRA SANITARIA 35G RUBI FLORAL'
$data[23,6] = 176
$data[23,7] = 13.2
$data[23,8] = 8.65
$data[24,0] = '2025-07-17'
$data[24,1] = 100
$data[24,2] = 'TEL TELECOMUNICACOES LTDA.'
$data[24,3] = 54443314
$data[24,4] = 11936640
$data[24,5] = 'LIMPADOR VEJA MULTIUSO GOLD 500ML'
$data[24,6] = 3383
$data[24,7] = 10.62
$data[24,8] = 14.65
$data[25,0] = '2025-07-21'
$data[25,1] = 200
$data[25,2] = 'JURUA ESTALEIROS E NAVEGACAO LTDA'
$data[25,3] = 54503121
$data[25,4] = 13996941
$data[25,5] = 'SACO DE LIXO 30L REFORCADO PACOTINHO C/10 UND FORTE MAX'
$data[25,6] = 91
$data[25,7] = 13.49
$data[25,8] = 22.27
$data[26,0] = '2025-07-21'
$data[26,1] = 240
$data[26,2] = 'JURUA ESTALEIROS E NAVEGACAO LTDA'
$data[26,3] = 54563773
$data[26,4] = 16933123
$data[26,5] = 'SABAO EM PO LAVAGEM PERFEITA OMO  - 400G'
$data[26,6] = 67
$data[26,7] = 37.23
$data[26,8] = 64.78
$data[27,0] = '2025-07-22'
$data[27,1] = 9
$data[27,2] = 'CARITAS ARQUIDIOCESANA DE MANAUS'
$data[27,3] = 54625975
$data[27,4] = 17059594
$data[27,5] = 'SABAO EM PO ESPUMIL 4KG'
$data[27,6] = 1
$data[27,7] = 1.94
$data[27,8] = 1.87
$data[28,0] = '2025-07-22'
$data[28,1] = 22
$data[28,2] = 'CONDOMINIO CRISTAL TOWER'
$data[28,3] = 54625981
$data[28,4] = 11939645
$data[28,5] = 'SACO DE LIXO 50L PRETO COMUM - PCT C/100 UND'
$data[28,6] = 154
$data[28,7] = 4.3
$data[28,8] = 4.78
$data[29,0] = '2025-07-23'
$data[29,1] = 400
$data[29,2] = 'V V REFEICOES LTDA'
$data[29,3] = 54692772
$data[29,4] = 32130390
$data[29,5] = 'ESPONJA MULTIUSO JEITOSA'
$data[29,6] = 11885
$data[29,7] = 21.36
$data[29,8] = 51.13

$ws.Rows(32).Delete()

$ws.Range("A2:I31").Value = $data
